$wb = $excel.ActiveWorkbook

# Rename "Ecommerce" sheet to "ECommerce" (capitalization fix)
$wsCommerce = $wb.Worksheets.Item("Ecommerce")
$wsCommerce.Name = "ECommerce"

# Make ECommerce the active sheet/tab (this also clears tabSelected on
# the previously-active REG sheet), with B23 selected (was B34)
$wsCommerce.Activate()
$wsCommerce.Range("B23").Select()
